$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column L header: "end time"
$ws.Range("L1").Value = "end time"

# Row 10: walltime correction (8-core Macpro)
$ws.Range("H10").Value = 13

# Row 11: add walltime + running total (time-of-day) for 8-core Macpro sim
$ws.Range("H11").Value = 15
$ws.Range("J11").Value = 1.6595833333333334
$ws.Range("J11").NumberFormat = "h:mm AM/PM"

# Row 12: add walltime, final run note, running total + end time for 8-core Macpro sim
$ws.Range("H12").Value = 18
$ws.Range("I12").Value = "total = 2.5 days"
$ws.Range("J12").Value = 0.35972222222222222
$ws.Range("J12").NumberFormat = "h:mm AM/PM"
$ws.Range("L12").Value = 0.35972222222222222
$ws.Range("L12").NumberFormat = "h:mm AM/PM"

# Row 6: walltime correction (6-core Macpro)
$ws.Range("H6").Value = 11

# Row 7: add walltime + running total (time-of-day) for 6-core Macpro sim
$ws.Range("H7").Value = 13
$ws.Range("J7").Value = 0.65902777777777777
$ws.Range("J7").NumberFormat = "h:mm AM/PM"

# Row 8: final run note + running total (time-of-day) for 6-core Macpro sim
$ws.Range("I8").Value = "total = 2.3 days"
$ws.Range("J8").Value = 0.20069444444444443
$ws.Range("J8").NumberFormat = "h:mm AM/PM"

# Row 13: additional simulation summary
$ws.Range("E13").Value = 1200
$ws.Range("F13").Value = 1200
$ws.Range("H13").Formula = "=SUM(H9:H12)"
$ws.Range("I13").Formula = "=(G13*H13)/G12"
$ws.Range("J13").Formula = "=I13/24"

# Row 14: additional simulation summary
$ws.Range("E14").Value = 1200
$ws.Range("F14").Value = 4800
$ws.Range("H14").Formula = "=SUM(H10:H13)"
$ws.Range("I14").Formula = "=(G14*H14)/G13"
$ws.Range("J14").Formula = "=I14/24"

# Restore selection to where the author left off
$ws.Range("H18").Select() | Out-Null
